$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.858.43'
$ws.Range('E2').Value = '  +1.76%  '
$ws.Range('D3').Value = '3.498.33'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'600.51"
$ws.Range('E5').Value = '  +2.45%  '
$ws.Range('D6').Value = "'171.15"
$ws.Range('E6').Value = '  +2.24%  '
$ws.Range('D7').Value = "'0.607"
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('D8').Value = '3.493.48'
$ws.Range('E8').Value = '  +0.68%  '
$ws.Range('D9').Value = "'1.00"
$ws.Range('E10').Value = '  +0.77%  '
$ws.Range('D11').Value = "'7.26"
$ws.Range('E11').Value = '  +6.78%  '
$ws.Range('D12').Value = "'0.578"
$ws.Range('E12').Value = '  +1.12%  '
$ws.Range('D13').Value = "'45.85"
$ws.Range('E13').Value = '  -1.19%  '
$ws.Range('D14').Value = "'0.0000273"
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').Value = '4.066.60'
$ws.Range('E15').Value = '  +0.88%  '
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = "'601.22"
$ws.Range('E17').Value = '  -2.46%  '
$ws.Range('D18').Value = '3.497.29'
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('D19').Value = '69.927.28'
$ws.Range('E19').Value = '  +1.83%  '
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('E21').Value = '  -1.00%  '
$ws.Range('D22').Value = "'0.864"
$ws.Range('E22').Value = '  -0.97%  '
$ws.Range('D23').Value = "'9.13"
$ws.Range('E23').Value = '  -17.13%  '
$ws.Range('D24').Value = "'15.41"
$ws.Range('E24').Value = '  -2.09%  '
$ws.Range('D25').Value = "'95.06"
$ws.Range('E25').Value = '  -0.66%  '
$ws.Range('E26').Value = '  -1.97%  '
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('E28').Value = '  -1.39%  '
$ws.Range('D29').Value = "'33.72"
$ws.Range('E29').Value = '  +3.01%  '
$ws.Range('D30').Value = "'8.91"
$ws.Range('E30').Value = '  -2.17%  '
$ws.Range('D31').Value = "'696.33"
$ws.Range('E31').Value = '  +20.37%  '
$ws.Range('D32').Value = "'3.00"
$ws.Range('E32').Value = '  -2.62%  '
$ws.Range('D33').Value = "'8.05"
$ws.Range('E33').Value = '  -4.10%  '
$ws.Range('D34').Value = "'6.85"
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('E35').Value = '  -1.87%  '
$ws.Range('E36').Value = '  -1.94%  '
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('E38').Value = '  -0.49%  '
$ws.Range('E39').Value = '  +7.70%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = "'56.69"
$ws.Range('E40').Value = '  +0.14%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = "'0.998"
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('E42').Value = '  +3.86%  '
$ws.Range('D43').Value = '3.313.30'
$ws.Range('E43').Value = '  -2.32%  '
$ws.Range('D44').Value = "'0.311"
$ws.Range('E44').Value = '  -3.34%  '
$ws.Range('D45').Value = "'2.90"
$ws.Range('E45').Value = '  +3.65%  '
$ws.Range('D46').Value = "'32.04"
$ws.Range('E46').Value = '  -1.68%  '
$ws.Range('D47').Value = '0.0₃0684'
$ws.Range('E47').Value = '  -0.96%  '
$ws.Range('D48').Value = "'2.53"
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('E49').Value = '  +0.55%  '
$ws.Range('D50').Value = "'132.82"
$ws.Range('E50').Value = '  +0.07%  '
$ws.Range('E51').Value = '  -0.01%  '
